$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Cells.Item(11, 3).Value = 17.14285714285714
$ws1.Cells.Item(11, 4).Value = 6
$ws1.Cells.Item(14, 3).Value = 42.85714285714285
$ws1.Cells.Item(14, 4).Value = 48
$ws1.Cells.Item(14, 5).Value = 112
$ws1.Cells.Item(15, 3).Value = 33.33333333333333
$ws1.Cells.Item(15, 4).Value = 24
$ws1.Cells.Item(19, 3).Value = 16.21621621621622
$ws1.Cells.Item(19, 4).Value = 6
$ws1.Cells.Item(19, 5).Value = 37
$ws1.Cells.Item(23, 3).Value = 5.88235294117647
$ws1.Cells.Item(23, 4).Value = 1
$ws1.Cells.Item(23, 5).Value = 17
$ws1.Cells.Item(27, 3).Value = 32.60869565217391
$ws1.Cells.Item(27, 4).Value = 15
$ws1.Cells.Item(27, 5).Value = 46
$ws1.Cells.Item(33, 3).Value = 35.35353535353536
$ws1.Cells.Item(33, 4).Value = 35
$ws1.Cells.Item(33, 5).Value = 99
$ws1.Cells.Item(36, 3).Value = 36.98630136986301
$ws1.Cells.Item(36, 4).Value = 27
$ws1.Cells.Item(40, 3).Value = 34.73684210526316
$ws1.Cells.Item(40, 4).Value = 33
$ws1.Cells.Item(44, 3).Value = 23.28767123287671
$ws1.Cells.Item(44, 4).Value = 17
$ws1.Cells.Item(44, 5).Value = 73
$ws1.Cells.Item(47, 3).Value = 51.61290322580645
$ws1.Cells.Item(47, 4).Value = 32
$ws1.Cells.Item(49, 3).Value = 8.474576271186439
$ws1.Cells.Item(49, 4).Value = 5
$ws1.Cells.Item(53, 3).Value = 46.15384615384615
$ws1.Cells.Item(53, 4).Value = 18
$ws1.Cells.Item(61, 3).Value = 47.82608695652174
$ws1.Cells.Item(61, 5).Value = 115
$ws1.Cells.Item(65, 3).Value = 35.9375
$ws1.Cells.Item(65, 5).Value = 64
$ws1.Cells.Item(66, 3).Value = 19.17808219178082
$ws1.Cells.Item(66, 4).Value = 14
$ws1.Cells.Item(69, 3).Value = 45.58823529411764
$ws1.Cells.Item(69, 4).Value = 31
$ws1.Cells.Item(69, 5).Value = 68
$ws1.Cells.Item(73, 3).Value = 15.95744680851064
$ws1.Cells.Item(73, 4).Value = 15
$ws1.Cells.Item(73, 5).Value = 94
$ws1.Cells.Item(76, 3).Value = 27.3972602739726
$ws1.Cells.Item(76, 5).Value = 73
$ws1.Cells.Item(93, 3).Value = 45.16129032258064
$ws1.Cells.Item(93, 4).Value = 28
$ws1.Cells.Item(97, 3).Value = 62.7906976744186
$ws1.Cells.Item(97, 4).Value = 27
$ws1.Cells.Item(97, 5).Value = 43
$ws1.Cells.Item(99, 3).Value = 13.95348837209302
$ws1.Cells.Item(99, 4).Value = 12
$ws1.Cells.Item(99, 5).Value = 86
$ws1.Cells.Item(101, 3).Value = 7.407407407407407
$ws1.Cells.Item(101, 4).Value = 4
$ws1.Cells.Item(102, 3).Value = 4.761904761904762
$ws1.Cells.Item(102, 4).Value = 2
$ws1.Cells.Item(104, 3).Value = 3.875968992248062
$ws1.Cells.Item(104, 4).Value = 5
$ws1.Cells.Item(106, 3).Value = 2.857142857142857
$ws1.Cells.Item(106, 4).Value = 1
$ws1.Cells.Item(108, 3).Value = 4.166666666666666
$ws1.Cells.Item(108, 4).Value = 2
$ws1.Cells.Item(111, 3).Value = 3.305785123966942
$ws1.Cells.Item(111, 4).Value = 8
$ws1.Cells.Item(111, 5).Value = 242
$ws1.Cells.Item(112, 3).Value = 6.185567010309279
$ws1.Cells.Item(112, 4).Value = 6
$ws1.Cells.Item(112, 5).Value = 97
$ws1.Cells.Item(123, 3).Value = 4.347826086956522
$ws1.Cells.Item(123, 4).Value = 4
$ws1.Cells.Item(126, 3).Value = 6.666666666666667
$ws1.Cells.Item(126, 4).Value = 7
$ws1.Cells.Item(130, 3).Value = 3.550295857988166
$ws1.Cells.Item(130, 4).Value = 6
$ws1.Cells.Item(130, 5).Value = 169
$ws1.Cells.Item(132, 3).Value = 5.208333333333334
$ws1.Cells.Item(132, 4).Value = 10
$ws1.Cells.Item(133, 3).Value = 6.428571428571428
$ws1.Cells.Item(133, 4).Value = 9
$ws1.Cells.Item(138, 3).Value = 9.090909090909092
$ws1.Cells.Item(138, 4).Value = 5
$ws1.Cells.Item(140, 3).Value = 2.173913043478261
$ws1.Cells.Item(140, 4).Value = 1
$ws1.Cells.Item(144, 3).Value = 5.376344086021505
$ws1.Cells.Item(144, 4).Value = 5
$ws1.Cells.Item(146, 3).Value = 1.851851851851852
$ws1.Cells.Item(146, 4).Value = 1
$ws1.Cells.Item(148, 3).Value = 4.597701149425287
$ws1.Cells.Item(148, 4).Value = 4
$ws1.Cells.Item(148, 5).Value = 87
$ws1.Cells.Item(154, 3).Value = 0
$ws1.Cells.Item(154, 4).Value = 0
$ws1.Cells.Item(154, 5).Value = 26
$ws1.Cells.Item(159, 3).Value = 5.88235294117647
$ws1.Cells.Item(159, 4).Value = 6
$ws1.Cells.Item(161, 3).Value = 6.896551724137931
$ws1.Cells.Item(161, 4).Value = 12
$ws1.Cells.Item(162, 3).Value = 4.166666666666666
$ws1.Cells.Item(162, 4).Value = 4
$ws1.Cells.Item(163, 3).Value = 11.57894736842105
$ws1.Cells.Item(163, 4).Value = 11
$ws1.Cells.Item(165, 3).Value = 5
$ws1.Cells.Item(165, 4).Value = 3
$ws1.Cells.Item(166, 3).Value = 3.149606299212598
$ws1.Cells.Item(166, 5).Value = 127
$ws1.Cells.Item(168, 3).Value = 6.358381502890173
$ws1.Cells.Item(168, 4).Value = 11
$ws1.Cells.Item(168, 5).Value = 173
$ws1.Cells.Item(170, 3).Value = 7.446808510638298
$ws1.Cells.Item(170, 4).Value = 7
$ws1.Cells.Item(171, 3).Value = 6.024096385542169
$ws1.Cells.Item(171, 4).Value = 5
$ws1.Cells.Item(173, 3).Value = 8.988764044943821
$ws1.Cells.Item(173, 5).Value = 89
$ws1.Cells.Item(177, 3).Value = 4.545454545454546
$ws1.Cells.Item(177, 4).Value = 6
$ws1.Cells.Item(182, 3).Value = 4.294478527607362
$ws1.Cells.Item(182, 4).Value = 7
$ws1.Cells.Item(183, 3).Value = 3.260869565217391
$ws1.Cells.Item(183, 4).Value = 3
$ws1.Cells.Item(190, 3).Value = 8.130081300813007
$ws1.Cells.Item(190, 4).Value = 10
$ws1.Cells.Item(190, 5).Value = 123
$ws1.Cells.Item(191, 3).Value = 3.846153846153846
$ws1.Cells.Item(191, 4).Value = 3
$ws1.Cells.Item(196, 3).Value = 5.263157894736842
$ws1.Cells.Item(196, 4).Value = 5
$ws1.Cells.Item(201, 3).Value = 3.597122302158273
$ws1.Cells.Item(201, 4).Value = 5
$ws1.Cells.Item(204, 3).Value = 6.666666666666667
$ws1.Cells.Item(204, 4).Value = 2
$ws1.Cells.Item(205, 3).Value = 3.571428571428571
$ws1.Cells.Item(205, 4).Value = 2
$ws1.Cells.Item(206, 3).Value = 3.03030303030303
$ws1.Cells.Item(206, 4).Value = 2
$ws1.Cells.Item(208, 3).Value = 6.967213114754098
$ws1.Cells.Item(208, 4).Value = 17
$ws1.Cells.Item(208, 5).Value = 244
$ws1.Cells.Item(209, 3).Value = 2.884615384615385
$ws1.Cells.Item(209, 4).Value = 3
$ws1.Cells.Item(214, 3).Value = 3.703703703703703
$ws1.Cells.Item(214, 5).Value = 54
$ws1.Cells.Item(216, 3).Value = 5.825242718446602
$ws1.Cells.Item(216, 4).Value = 6
$ws1.Cells.Item(222, 3).Value = 2.816901408450704
$ws1.Cells.Item(222, 4).Value = 2
$ws1.Cells.Item(223, 3).Value = 4.132231404958678
$ws1.Cells.Item(223, 4).Value = 5
$ws1.Cells.Item(226, 3).Value = 7.407407407407407
$ws1.Cells.Item(226, 4).Value = 8
$ws1.Cells.Item(227, 3).Value = 2.352941176470588
$ws1.Cells.Item(227, 4).Value = 4
$ws1.Cells.Item(233, 3).Value = 3.448275862068965
$ws1.Cells.Item(233, 4).Value = 3
$ws1.Cells.Item(234, 3).Value = 3.636363636363636
$ws1.Cells.Item(234, 4).Value = 6
$ws1.Cells.Item(238, 3).Value = 6.382978723404255
$ws1.Cells.Item(238, 4).Value = 6
$ws1.Cells.Item(241, 3).Value = 0
$ws1.Cells.Item(241, 4).Value = 0
$ws1.Cells.Item(243, 3).Value = 5.172413793103448
$ws1.Cells.Item(243, 4).Value = 3
$ws1.Cells.Item(251, 3).Value = 0
$ws1.Cells.Item(251, 4).Value = 0
$ws1.Cells.Item(251, 5).Value = 28
$ws1.Cells.Item(252, 3).Value = 0.8695652173913043
$ws1.Cells.Item(252, 4).Value = 1
$ws1.Cells.Item(256, 3).Value = 2.857142857142857
$ws1.Cells.Item(256, 4).Value = 3
$ws1.Cells.Item(258, 3).Value = 6.074766355140187
$ws1.Cells.Item(258, 4).Value = 13
$ws1.Cells.Item(259, 3).Value = 5.660377358490567
$ws1.Cells.Item(259, 4).Value = 6
$ws1.Cells.Item(260, 3).Value = 2.97029702970297
$ws1.Cells.Item(260, 4).Value = 3
$ws1.Cells.Item(262, 3).Value = 3.508771929824561
$ws1.Cells.Item(262, 4).Value = 2
$ws1.Cells.Item(263, 3).Value = 5.109489051094891
$ws1.Cells.Item(263, 4).Value = 7
$ws1.Cells.Item(263, 5).Value = 137
$ws1.Cells.Item(265, 3).Value = 4.790419161676647
$ws1.Cells.Item(265, 4).Value = 8
$ws1.Cells.Item(265, 5).Value = 167
$ws1.Cells.Item(267, 3).Value = 3.361344537815126
$ws1.Cells.Item(267, 4).Value = 4
$ws1.Cells.Item(270, 3).Value = 6.60377358490566
$ws1.Cells.Item(270, 5).Value = 106
$ws1.Cells.Item(274, 3).Value = 3.597122302158273
$ws1.Cells.Item(274, 4).Value = 5
$ws1.Cells.Item(275, 3).Value = 0
$ws1.Cells.Item(275, 4).Value = 0
$ws1.Cells.Item(279, 3).Value = 6.707317073170732
$ws1.Cells.Item(279, 4).Value = 11
$ws1.Cells.Item(285, 3).Value = 9.45945945945946
$ws1.Cells.Item(285, 4).Value = 7
$ws1.Cells.Item(285, 5).Value = 74
$ws1.Cells.Item(287, 3).Value = 6.923076923076923
$ws1.Cells.Item(287, 4).Value = 9
$ws1.Cells.Item(287, 5).Value = 130
$ws1.Cells.Item(295, 3).Value = 0
$ws1.Cells.Item(295, 4).Value = 0
$ws1.Cells.Item(300, 3).Value = 2.083333333333333
$ws1.Cells.Item(300, 4).Value = 1
$ws1.Cells.Item(302, 3).Value = 3.076923076923077
$ws1.Cells.Item(302, 4).Value = 2
$ws1.Cells.Item(305, 3).Value = 1.094890510948905
$ws1.Cells.Item(305, 4).Value = 3
$ws1.Cells.Item(305, 5).Value = 274
$ws1.Cells.Item(306, 3).Value = 0
$ws1.Cells.Item(306, 4).Value = 0
$ws1.Cells.Item(313, 3).Value = 2.521008403361344
$ws1.Cells.Item(313, 4).Value = 3
$ws1.Cells.Item(316, 3).Value = 2.02020202020202
$ws1.Cells.Item(316, 4).Value = 2
$ws1.Cells.Item(323, 3).Value = 0.7462686567164178
$ws1.Cells.Item(323, 4).Value = 1
$ws1.Cells.Item(324, 3).Value = 1.515151515151515
$ws1.Cells.Item(324, 4).Value = 3
$ws1.Cells.Item(327, 3).Value = 0.6329113924050633
$ws1.Cells.Item(327, 4).Value = 1
$ws1.Cells.Item(328, 3).Value = 0.6711409395973155
$ws1.Cells.Item(328, 4).Value = 1
$ws1.Cells.Item(330, 3).Value = 1.08695652173913
$ws1.Cells.Item(330, 4).Value = 1
$ws1.Cells.Item(331, 3).Value = 1.075268817204301
$ws1.Cells.Item(331, 4).Value = 2
$ws1.Cells.Item(337, 3).Value = 2.285714285714286
$ws1.Cells.Item(337, 4).Value = 4
$ws1.Cells.Item(338, 3).Value = 0.9090909090909091
$ws1.Cells.Item(338, 4).Value = 1
$ws1.Cells.Item(340, 3).Value = 0
$ws1.Cells.Item(340, 4).Value = 0
$ws1.Cells.Item(349, 3).Value = 0.7142857142857143
$ws1.Cells.Item(349, 4).Value = 1
$ws1.Cells.Item(352, 3).Value = 2.811244979919679
$ws1.Cells.Item(352, 4).Value = 7
$ws1.Cells.Item(356, 3).Value = 0.8333333333333334
$ws1.Cells.Item(356, 5).Value = 120
$ws1.Cells.Item(357, 3).Value = 0
$ws1.Cells.Item(357, 4).Value = 0
$ws1.Cells.Item(359, 3).Value = 0
$ws1.Cells.Item(359, 4).Value = 0
$ws1.Cells.Item(360, 3).Value = 1.342281879194631
$ws1.Cells.Item(360, 4).Value = 2
$ws1.Cells.Item(360, 5).Value = 149
$ws1.Cells.Item(363, 3).Value = 3.389830508474576
$ws1.Cells.Item(363, 4).Value = 2
$ws1.Cells.Item(365, 3).Value = 0.9900990099009901
$ws1.Cells.Item(365, 4).Value = 1
$ws1.Cells.Item(367, 3).Value = 2.439024390243902
$ws1.Cells.Item(367, 5).Value = 123
$ws1.Cells.Item(377, 3).Value = 0.9345794392523363
$ws1.Cells.Item(377, 4).Value = 1
$ws1.Cells.Item(380, 3).Value = 2.409638554216868
$ws1.Cells.Item(380, 4).Value = 2
$ws1.Cells.Item(384, 3).Value = 4.444444444444445
$ws1.Cells.Item(384, 4).Value = 6
$ws1.Cells.Item(392, 3).Value = 50
$ws1.Cells.Item(392, 4).Value = 10
$ws1.Cells.Item(395, 3).Value = 60.60606060606061
$ws1.Cells.Item(395, 4).Value = 20
$ws1.Cells.Item(403, 3).Value = 48.48484848484848
$ws1.Cells.Item(403, 4).Value = 16
$ws1.Cells.Item(406, 3).Value = 33.33333333333333
$ws1.Cells.Item(406, 4).Value = 9
$ws1.Cells.Item(408, 3).Value = 69.23076923076923
$ws1.Cells.Item(408, 4).Value = 9
$ws1.Cells.Item(421, 3).Value = 38.33333333333334
$ws1.Cells.Item(421, 4).Value = 23
$ws1.Cells.Item(423, 3).Value = 31.57894736842105
$ws1.Cells.Item(423, 4).Value = 18
$ws1.Cells.Item(434, 3).Value = 29.41176470588236
$ws1.Cells.Item(434, 4).Value = 15
$ws1.Cells.Item(450, 3).Value = 50
$ws1.Cells.Item(450, 4).Value = 18
$ws1.Cells.Item(452, 3).Value = 49.12280701754386
$ws1.Cells.Item(452, 4).Value = 28
$ws1.Cells.Item(457, 3).Value = 50
$ws1.Cells.Item(457, 4).Value = 22
$ws1.Cells.Item(457, 5).Value = 44
$ws1.Cells.Item(459, 3).Value = 51.85185185185185
$ws1.Cells.Item(459, 4).Value = 28
$ws1.Cells.Item(459, 5).Value = 54
$ws1.Cells.Item(461, 3).Value = 58.62068965517241
$ws1.Cells.Item(461, 5).Value = 29
$ws1.Cells.Item(463, 3).Value = 30
$ws1.Cells.Item(463, 4).Value = 6
$ws1.Cells.Item(464, 3).Value = 50
$ws1.Cells.Item(464, 4).Value = 14
$ws1.Cells.Item(464, 5).Value = 28
$ws1.Cells.Item(476, 3).Value = 38.09523809523809
$ws1.Cells.Item(476, 4).Value = 8
$ws1.Cells.Item(478, 3).Value = 60
$ws1.Cells.Item(478, 4).Value = 18
$ws1.Cells.Item(478, 5).Value = 30
$ws1.Cells.Item(481, 3).Value = 40.42553191489361
$ws1.Cells.Item(481, 4).Value = 19
$ws1.Cells.Item(487, 3).Value = 30
$ws1.Cells.Item(487, 4).Value = 15
$ws1.Cells.Item(496, 3).Value = 11.53846153846154
$ws1.Cells.Item(496, 4).Value = 3
$ws1.Cells.Item(500, 3).Value = 18.18181818181818
$ws1.Cells.Item(500, 4).Value = 8
$ws1.Cells.Item(504, 3).Value = 42.30769230769231
$ws1.Cells.Item(504, 4).Value = 11
$ws1.Cells.Item(518, 3).Value = 14.54545454545454
$ws1.Cells.Item(518, 4).Value = 16
$ws1.Cells.Item(520, 3).Value = 13.26530612244898
$ws1.Cells.Item(520, 4).Value = 13
$ws1.Cells.Item(521, 3).Value = 15.11627906976744
$ws1.Cells.Item(521, 4).Value = 13
$ws1.Cells.Item(528, 3).Value = 21.42857142857143
$ws1.Cells.Item(528, 4).Value = 6
$ws1.Cells.Item(531, 3).Value = 21.05263157894737
$ws1.Cells.Item(531, 4).Value = 16
$ws1.Cells.Item(532, 3).Value = 23.21428571428572
$ws1.Cells.Item(532, 4).Value = 13
$ws1.Cells.Item(552, 3).Value = 5
$ws1.Cells.Item(552, 4).Value = 1
$ws1.Cells.Item(553, 3).Value = 17.07317073170732
$ws1.Cells.Item(553, 4).Value = 7
$ws1.Cells.Item(554, 3).Value = 20.37037037037037
$ws1.Cells.Item(554, 4).Value = 11
$ws1.Cells.Item(554, 5).Value = 54
$ws1.Cells.Item(561, 3).Value = 31.25
$ws1.Cells.Item(561, 4).Value = 15
$ws1.Cells.Item(561, 5).Value = 48
$ws1.Cells.Item(571, 3).Value = 4.615384615384616
$ws1.Cells.Item(571, 4).Value = 3
$ws1.Cells.Item(578, 3).Value = 18.46153846153846
$ws1.Cells.Item(578, 4).Value = 12
$ws1.Cells.Item(582, 3).Value = 22.22222222222222
$ws1.Cells.Item(582, 4).Value = 12
$ws1.Cells.Item(586, 3).Value = 0
$ws1.Cells.Item(586, 4).Value = 0
$ws1.Cells.Item(589, 3).Value = 0.6896551724137931
$ws1.Cells.Item(589, 4).Value = 1
$ws1.Cells.Item(593, 3).Value = 1.538461538461539
$ws1.Cells.Item(593, 4).Value = 1
$ws1.Cells.Item(594, 3).Value = 0
$ws1.Cells.Item(594, 4).Value = 0
$ws1.Cells.Item(596, 3).Value = 0.7352941176470588
$ws1.Cells.Item(596, 4).Value = 2
$ws1.Cells.Item(596, 5).Value = 272
$ws1.Cells.Item(597, 3).Value = 0
$ws1.Cells.Item(597, 4).Value = 0
$ws1.Cells.Item(611, 3).Value = 0.6896551724137931
$ws1.Cells.Item(611, 4).Value = 1
$ws1.Cells.Item(615, 3).Value = 1.507537688442211
$ws1.Cells.Item(615, 4).Value = 3
$ws1.Cells.Item(618, 3).Value = 0
$ws1.Cells.Item(618, 4).Value = 0
$ws1.Cells.Item(622, 3).Value = 1.081081081081081
$ws1.Cells.Item(622, 4).Value = 2
$ws1.Cells.Item(629, 3).Value = 0
$ws1.Cells.Item(629, 4).Value = 0
$ws1.Cells.Item(631, 3).Value = 0
$ws1.Cells.Item(631, 4).Value = 0
$ws1.Cells.Item(646, 3).Value = 0.8771929824561403
$ws1.Cells.Item(646, 4).Value = 2
$ws1.Cells.Item(647, 3).Value = 0.8333333333333334
$ws1.Cells.Item(647, 4).Value = 1
$ws1.Cells.Item(647, 5).Value = 120
$ws1.Cells.Item(648, 3).Value = 0
$ws1.Cells.Item(648, 4).Value = 0
$ws1.Cells.Item(650, 3).Value = 0
$ws1.Cells.Item(650, 4).Value = 0
$ws1.Cells.Item(651, 3).Value = 1.351351351351351
$ws1.Cells.Item(651, 5).Value = 148
$ws1.Cells.Item(653, 3).Value = 1.522842639593909
$ws1.Cells.Item(653, 4).Value = 3
$ws1.Cells.Item(653, 5).Value = 197
$ws1.Cells.Item(658, 3).Value = 2.419354838709677
$ws1.Cells.Item(658, 4).Value = 3
$ws1.Cells.Item(658, 5).Value = 124
$ws1.Cells.Item(662, 3).Value = 2.068965517241379
$ws1.Cells.Item(662, 4).Value = 3
$ws1.Cells.Item(668, 3).Value = 1.886792452830189
$ws1.Cells.Item(668, 4).Value = 2
$ws1.Cells.Item(675, 3).Value = 5.970149253731343
$ws1.Cells.Item(675, 4).Value = 8
$ws1.Cells.Item(675, 5).Value = 134
$ws1.Cells.Item(683, 3).Value = 0
$ws1.Cells.Item(683, 4).Value = 0
$ws1.Cells.Item(686, 3).Value = 0.6896551724137931
$ws1.Cells.Item(686, 4).Value = 1
$ws1.Cells.Item(688, 3).Value = 2.083333333333333
$ws1.Cells.Item(688, 4).Value = 1
$ws1.Cells.Item(690, 3).Value = 3.076923076923077
$ws1.Cells.Item(690, 4).Value = 2
$ws1.Cells.Item(691, 3).Value = 0
$ws1.Cells.Item(691, 4).Value = 0
$ws1.Cells.Item(693, 3).Value = 1.454545454545455
$ws1.Cells.Item(693, 4).Value = 4
$ws1.Cells.Item(693, 5).Value = 275
$ws1.Cells.Item(694, 3).Value = 0
$ws1.Cells.Item(694, 4).Value = 0
$ws1.Cells.Item(701, 3).Value = 2.521008403361344
$ws1.Cells.Item(701, 4).Value = 3
$ws1.Cells.Item(704, 3).Value = 3
$ws1.Cells.Item(704, 4).Value = 3
$ws1.Cells.Item(711, 3).Value = 0.7462686567164178
$ws1.Cells.Item(711, 4).Value = 1
$ws1.Cells.Item(712, 3).Value = 2.010050251256281
$ws1.Cells.Item(712, 4).Value = 4
$ws1.Cells.Item(715, 3).Value = 0.6329113924050633
$ws1.Cells.Item(715, 4).Value = 1
$ws1.Cells.Item(716, 3).Value = 0.6711409395973155
$ws1.Cells.Item(716, 4).Value = 1
$ws1.Cells.Item(718, 3).Value = 2.173913043478261
$ws1.Cells.Item(718, 4).Value = 2
$ws1.Cells.Item(719, 3).Value = 1.612903225806452
$ws1.Cells.Item(719, 4).Value = 3
$ws1.Cells.Item(725, 3).Value = 2.857142857142857
$ws1.Cells.Item(725, 4).Value = 5
$ws1.Cells.Item(726, 3).Value = 0.9090909090909091
$ws1.Cells.Item(726, 4).Value = 1
$ws1.Cells.Item(728, 3).Value = 0
$ws1.Cells.Item(728, 4).Value = 0
$ws1.Cells.Item(737, 3).Value = 0.7142857142857143
$ws1.Cells.Item(737, 4).Value = 1
$ws1.Cells.Item(740, 3).Value = 3.614457831325301
$ws1.Cells.Item(740, 4).Value = 9
$ws1.Cells.Item(743, 3).Value = 1.739130434782609
$ws1.Cells.Item(743, 4).Value = 4
$ws1.Cells.Item(744, 3).Value = 0.8333333333333334
$ws1.Cells.Item(744, 4).Value = 1
$ws1.Cells.Item(744, 5).Value = 120
$ws1.Cells.Item(745, 3).Value = 0
$ws1.Cells.Item(745, 4).Value = 0
$ws1.Cells.Item(747, 3).Value = 0
$ws1.Cells.Item(747, 4).Value = 0
$ws1.Cells.Item(748, 3).Value = 1.342281879194631
$ws1.Cells.Item(748, 4).Value = 2
$ws1.Cells.Item(748, 5).Value = 149
$ws1.Cells.Item(750, 3).Value = 2.02020202020202
$ws1.Cells.Item(750, 4).Value = 4
$ws1.Cells.Item(751, 3).Value = 5.084745762711865
$ws1.Cells.Item(751, 4).Value = 3
$ws1.Cells.Item(753, 3).Value = 0.9900990099009901
$ws1.Cells.Item(753, 4).Value = 1
$ws1.Cells.Item(755, 3).Value = 2.419354838709677
$ws1.Cells.Item(755, 4).Value = 3
$ws1.Cells.Item(755, 5).Value = 124
$ws1.Cells.Item(759, 3).Value = 2.054794520547945
$ws1.Cells.Item(759, 4).Value = 3
$ws1.Cells.Item(765, 3).Value = 1.869158878504673
$ws1.Cells.Item(765, 4).Value = 2
$ws1.Cells.Item(768, 3).Value = 2.409638554216868
$ws1.Cells.Item(768, 4).Value = 2
$ws1.Cells.Item(772, 3).Value = 6.666666666666667
$ws1.Cells.Item(772, 4).Value = 9
$ws1.Cells.Item(787, 3).Value = 1.538461538461539
$ws1.Cells.Item(787, 4).Value = 1
$ws1.Cells.Item(790, 3).Value = 0.3636363636363636
$ws1.Cells.Item(790, 4).Value = 1
$ws1.Cells.Item(790, 5).Value = 275
$ws1.Cells.Item(791, 3).Value = 0
$ws1.Cells.Item(791, 4).Value = 0
$ws1.Cells.Item(805, 3).Value = 0.684931506849315
$ws1.Cells.Item(805, 4).Value = 1
$ws1.Cells.Item(809, 3).Value = 1.005025125628141
$ws1.Cells.Item(809, 4).Value = 2
$ws1.Cells.Item(812, 3).Value = 0
$ws1.Cells.Item(812, 4).Value = 0
$ws1.Cells.Item(816, 3).Value = 0.5376344086021506
$ws1.Cells.Item(816, 4).Value = 1
$ws1.Cells.Item(823, 3).Value = 0
$ws1.Cells.Item(823, 4).Value = 0
$ws1.Cells.Item(825, 3).Value = 0
$ws1.Cells.Item(825, 4).Value = 0
$ws1.Cells.Item(841, 3).Value = 0.8333333333333334
$ws1.Cells.Item(841, 5).Value = 120
$ws1.Cells.Item(842, 3).Value = 0
$ws1.Cells.Item(842, 4).Value = 0
$ws1.Cells.Item(845, 3).Value = 1.342281879194631
$ws1.Cells.Item(845, 5).Value = 149
$ws1.Cells.Item(852, 3).Value = 2.419354838709677
$ws1.Cells.Item(852, 5).Value = 124
$ws1.Cells.Item(862, 3).Value = 0.9345794392523363
$ws1.Cells.Item(862, 4).Value = 1
$ws1.Cells.Item(869, 3).Value = 3.703703703703703
$ws1.Cells.Item(869, 4).Value = 5

$ws2.Cells.Item(2, 4).Value = 5.1
$ws2.Cells.Item(2, 5).Value = 41
$ws2.Cells.Item(2, 6).Value = 804
$ws2.Cells.Item(4, 4).Value = 50.74
$ws2.Cells.Item(4, 5).Value = 138
$ws2.Cells.Item(5, 4).Value = 2.42
$ws2.Cells.Item(5, 5).Value = 21
$ws2.Cells.Item(6, 4).Value = 52.72
$ws2.Cells.Item(6, 5).Value = 194
$ws2.Cells.Item(6, 6).Value = 368
$ws2.Cells.Item(7, 4).Value = 5.69
$ws2.Cells.Item(7, 5).Value = 48
$ws2.Cells.Item(7, 6).Value = 843
$ws2.Cells.Item(8, 4).Value = 3.1
$ws2.Cells.Item(8, 5).Value = 27
$ws2.Cells.Item(9, 4).Value = 1.95
$ws2.Cells.Item(9, 5).Value = 17
$ws2.Cells.Item(10, 4).Value = 2.65
$ws2.Cells.Item(10, 5).Value = 23
$ws2.Cells.Item(10, 6).Value = 868
$ws2.Cells.Item(11, 4).Value = 7.84
$ws2.Cells.Item(11, 5).Value = 28
$ws2.Cells.Item(12, 4).Value = 23.32
$ws2.Cells.Item(12, 5).Value = 45
$ws2.Cells.Item(15, 4).Value = 32.41
$ws2.Cells.Item(15, 5).Value = 82
$ws2.Cells.Item(15, 6).Value = 253
$ws2.Cells.Item(16, 4).Value = 1.69
$ws2.Cells.Item(16, 5).Value = 7
$ws2.Cells.Item(18, 4).Value = 0
$ws2.Cells.Item(18, 5).Value = 0
$ws2.Cells.Item(19, 4).Value = 0.23
$ws2.Cells.Item(19, 5).Value = 1
$ws2.Cells.Item(20, 4).Value = 8.44
$ws2.Cells.Item(20, 5).Value = 40
$ws2.Cells.Item(22, 4).Value = 52.76
$ws2.Cells.Item(22, 6).Value = 163
$ws2.Cells.Item(24, 4).Value = 20.15
$ws2.Cells.Item(24, 5).Value = 82
$ws2.Cells.Item(24, 6).Value = 407
$ws2.Cells.Item(25, 4).Value = 6.8
$ws2.Cells.Item(25, 6).Value = 544
$ws2.Cells.Item(29, 4).Value = 6.57
$ws2.Cells.Item(29, 5).Value = 27
$ws2.Cells.Item(29, 6).Value = 411
$ws2.Cells.Item(30, 4).Value = 15.22
$ws2.Cells.Item(30, 5).Value = 35
$ws2.Cells.Item(31, 4).Value = 46.88
$ws2.Cells.Item(31, 5).Value = 60
$ws2.Cells.Item(32, 4).Value = 0.89
$ws2.Cells.Item(32, 5).Value = 5
$ws2.Cells.Item(33, 4).Value = 26.75
$ws2.Cells.Item(33, 6).Value = 314
$ws2.Cells.Item(34, 4).Value = 4.02
$ws2.Cells.Item(34, 5).Value = 20
$ws2.Cells.Item(35, 4).Value = 1.07
$ws2.Cells.Item(35, 5).Value = 6
$ws2.Cells.Item(36, 4).Value = 0.89
$ws2.Cells.Item(36, 5).Value = 5
$ws2.Cells.Item(37, 4).Value = 1.07
$ws2.Cells.Item(37, 5).Value = 6
$ws2.Cells.Item(38, 4).Value = 6.84
$ws2.Cells.Item(38, 5).Value = 44
$ws2.Cells.Item(40, 4).Value = 49.31
$ws2.Cells.Item(40, 5).Value = 107
$ws2.Cells.Item(41, 4).Value = 1.63
$ws2.Cells.Item(41, 5).Value = 13
$ws2.Cells.Item(42, 4).Value = 38.29
$ws2.Cells.Item(42, 6).Value = 397
$ws2.Cells.Item(43, 4).Value = 4.41
$ws2.Cells.Item(43, 5).Value = 32
$ws2.Cells.Item(46, 4).Value = 1.64
$ws2.Cells.Item(46, 5).Value = 13
$ws2.Cells.Item(47, 4).Value = 4.14
$ws2.Cells.Item(47, 5).Value = 29
$ws2.Cells.Item(48, 4).Value = 24.86
$ws2.Cells.Item(48, 5).Value = 91
$ws2.Cells.Item(48, 6).Value = 366
$ws2.Cells.Item(49, 4).Value = 47.74
$ws2.Cells.Item(49, 5).Value = 116
$ws2.Cells.Item(50, 4).Value = 1.86
$ws2.Cells.Item(50, 5).Value = 16
$ws2.Cells.Item(50, 6).Value = 862
$ws2.Cells.Item(51, 4).Value = 34.68
$ws2.Cells.Item(51, 5).Value = 146
$ws2.Cells.Item(51, 6).Value = 421
$ws2.Cells.Item(52, 4).Value = 4.42
$ws2.Cells.Item(52, 5).Value = 34
$ws2.Cells.Item(53, 4).Value = 1.86
$ws2.Cells.Item(53, 5).Value = 16
$ws2.Cells.Item(53, 6).Value = 862
$ws2.Cells.Item(54, 4).Value = 1.51
$ws2.Cells.Item(54, 5).Value = 13
$ws2.Cells.Item(54, 6).Value = 862
$ws2.Cells.Item(55, 4).Value = 1.52
$ws2.Cells.Item(55, 5).Value = 13
$ws2.Cells.Item(55, 6).Value = 858
$ws2.Cells.Item(56, 4).Value = 6.79
$ws2.Cells.Item(56, 6).Value = 442
$ws2.Cells.Item(57, 4).Value = 20.09
$ws2.Cells.Item(57, 5).Value = 46
$ws2.Cells.Item(58, 4).Value = 29.56
$ws2.Cells.Item(58, 5).Value = 47
$ws2.Cells.Item(65, 4).Value = 4.73
$ws2.Cells.Item(65, 5).Value = 20
$ws2.Cells.Item(68, 4).Value = 1.08
$ws2.Cells.Item(68, 5).Value = 6
$ws2.Cells.Item(70, 4).Value = 3.65
$ws2.Cells.Item(70, 5).Value = 17
$ws2.Cells.Item(71, 4).Value = 1.08
$ws2.Cells.Item(71, 5).Value = 6
$ws2.Cells.Item(72, 4).Value = 0.72
$ws2.Cells.Item(72, 5).Value = 4
$ws2.Cells.Item(73, 4).Value = 0.72
$ws2.Cells.Item(73, 5).Value = 4
$ws2.Cells.Item(74, 4).Value = 6.06
$ws2.Cells.Item(74, 5).Value = 53
$ws2.Cells.Item(75, 4).Value = 17.17
$ws2.Cells.Item(75, 5).Value = 79
$ws2.Cells.Item(77, 4).Value = 0.27
$ws2.Cells.Item(77, 5).Value = 3
$ws2.Cells.Item(78, 4).Value = 20.98
$ws2.Cells.Item(78, 5).Value = 128
$ws2.Cells.Item(78, 6).Value = 610
$ws2.Cells.Item(79, 4).Value = 2.64
$ws2.Cells.Item(79, 5).Value = 25
$ws2.Cells.Item(79, 6).Value = 948
$ws2.Cells.Item(80, 4).Value = 0.53
$ws2.Cells.Item(80, 5).Value = 6
$ws2.Cells.Item(81, 4).Value = 0.09
$ws2.Cells.Item(81, 5).Value = 1
$ws2.Cells.Item(82, 4).Value = 0.36
$ws2.Cells.Item(82, 5).Value = 4
$ws2.Cells.Item(83, 4).Value = 4.98
$ws2.Cells.Item(83, 5).Value = 43
$ws2.Cells.Item(83, 6).Value = 863
$ws2.Cells.Item(84, 4).Value = 14.36
$ws2.Cells.Item(84, 5).Value = 79
$ws2.Cells.Item(85, 4).Value = 43.4
$ws2.Cells.Item(85, 5).Value = 115
$ws2.Cells.Item(86, 4).Value = 0.82
$ws2.Cells.Item(86, 5).Value = 9
$ws2.Cells.Item(87, 4).Value = 26.91
$ws2.Cells.Item(87, 5).Value = 155
$ws2.Cells.Item(87, 6).Value = 576
$ws2.Cells.Item(88, 4).Value = 4.12
$ws2.Cells.Item(88, 5).Value = 37
$ws2.Cells.Item(89, 4).Value = 1.09
$ws2.Cells.Item(89, 5).Value = 12
$ws2.Cells.Item(90, 4).Value = 0.64
$ws2.Cells.Item(90, 5).Value = 7
$ws2.Cells.Item(91, 4).Value = 0.91
$ws2.Cells.Item(91, 5).Value = 10
$ws2.Cells.Item(92, 4).Value = 7.88
$ws2.Cells.Item(92, 5).Value = 80
$ws2.Cells.Item(92, 6).Value = 1015
$ws2.Cells.Item(93, 4).Value = 21.6
$ws2.Cells.Item(93, 5).Value = 119
$ws2.Cells.Item(93, 6).Value = 551
$ws2.Cells.Item(94, 4).Value = 50.32
$ws2.Cells.Item(94, 5).Value = 156
$ws2.Cells.Item(94, 6).Value = 310
$ws2.Cells.Item(95, 4).Value = 0.98
$ws2.Cells.Item(95, 5).Value = 13
$ws2.Cells.Item(95, 6).Value = 1333
$ws2.Cells.Item(96, 4).Value = 26.86
$ws2.Cells.Item(96, 5).Value = 199
$ws2.Cells.Item(96, 6).Value = 741
$ws2.Cells.Item(97, 4).Value = 4.91
$ws2.Cells.Item(97, 5).Value = 54
$ws2.Cells.Item(98, 4).Value = 1.35
$ws2.Cells.Item(98, 5).Value = 18
$ws2.Cells.Item(98, 6).Value = 1336
$ws2.Cells.Item(99, 4).Value = 0.75
$ws2.Cells.Item(99, 5).Value = 10
$ws2.Cells.Item(99, 6).Value = 1336
$ws2.Cells.Item(100, 4).Value = 1.13
$ws2.Cells.Item(100, 5).Value = 15
$ws2.Cells.Item(100, 6).Value = 1329
$ws2.Cells.Item(101, 4).Value = 3.7
$ws2.Cells.Item(101, 5).Value = 26
$ws2.Cells.Item(101, 6).Value = 702
$ws2.Cells.Item(102, 4).Value = 14.07
$ws2.Cells.Item(102, 5).Value = 56
$ws2.Cells.Item(103, 4).Value = 42.93
$ws2.Cells.Item(103, 5).Value = 82
$ws2.Cells.Item(104, 4).Value = 1.11
$ws2.Cells.Item(104, 5).Value = 9
$ws2.Cells.Item(104, 6).Value = 808
$ws2.Cells.Item(105, 4).Value = 37.39
$ws2.Cells.Item(105, 5).Value = 126
$ws2.Cells.Item(105, 6).Value = 337
$ws2.Cells.Item(106, 4).Value = 6.28
$ws2.Cells.Item(106, 5).Value = 45
$ws2.Cells.Item(106, 6).Value = 716
$ws2.Cells.Item(107, 4).Value = 1.61
$ws2.Cells.Item(107, 5).Value = 13
$ws2.Cells.Item(107, 6).Value = 809
$ws2.Cells.Item(108, 4).Value = 0.74
$ws2.Cells.Item(108, 5).Value = 6
$ws2.Cells.Item(108, 6).Value = 809
$ws2.Cells.Item(109, 4).Value = 1.24
$ws2.Cells.Item(109, 5).Value = 10
$ws2.Cells.Item(109, 6).Value = 805

$ws3.Cells.Item(2, 2).Value = 6.06
$ws3.Cells.Item(2, 3).Value = 468
$ws3.Cells.Item(2, 4).Value = 7718
$ws3.Cells.Item(3, 2).Value = 19.26
$ws3.Cells.Item(3, 3).Value = 823
$ws3.Cells.Item(3, 4).Value = 4272
$ws3.Cells.Item(4, 2).Value = 45.31
$ws3.Cells.Item(4, 3).Value = 1120
$ws3.Cells.Item(4, 4).Value = 2472
$ws3.Cells.Item(5, 2).Value = 1.16
$ws3.Cells.Item(5, 3).Value = 114
$ws3.Cells.Item(5, 4).Value = 9794
$ws3.Cells.Item(6, 2).Value = 29
$ws3.Cells.Item(6, 3).Value = 1496
$ws3.Cells.Item(6, 4).Value = 5159
$ws3.Cells.Item(7, 2).Value = 4.55
$ws3.Cells.Item(7, 3).Value = 383
$ws3.Cells.Item(8, 2).Value = 1.5
$ws3.Cells.Item(8, 3).Value = 147
$ws3.Cells.Item(8, 4).Value = 9808
$ws3.Cells.Item(9, 2).Value = 0.85
$ws3.Cells.Item(9, 3).Value = 83
$ws3.Cells.Item(9, 4).Value = 9808
$ws3.Cells.Item(10, 2).Value = 1.19
$ws3.Cells.Item(10, 3).Value = 116
